$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header E5 text: "$ CHG " -> "$ CHG" (new shared string) ---
$ws.Cells.Item(5, 5).Value = "$ CHG"

# --- Row 6: replace formulas in E6/F6 with static values, update B6/C6/D6 ---
$ws.Cells.Item(6, 2).Value = 750
$ws.Cells.Item(6, 3).Value = 500
$ws.Cells.Item(6, 4).Value = 440
$ws.Cells.Item(6, 5).Value = 250
$ws.Cells.Item(6, 6).Value = 0.5
# G6 (0.1) is unchanged

# --- Add new rows 7-9, copying the B:F number formatting/styles from row 6 ---
$ws.Range("B6:F6").Copy()
$ws.Range("B7:F9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 7 values ---
$ws.Cells.Item(7, 2).Value = 500
$ws.Cells.Item(7, 3).Value = 640
$ws.Cells.Item(7, 4).Value = 470
$ws.Cells.Item(7, 5).Value = -140
$ws.Cells.Item(7, 6).Value = -0.21875

# --- Row 8 values ---
$ws.Cells.Item(8, 2).Value = 600
$ws.Cells.Item(8, 3).Value = 500
$ws.Cells.Item(8, 4).Value = 450
$ws.Cells.Item(8, 5).Value = 100
$ws.Cells.Item(8, 6).Value = 0.2

# --- Row 9 values ---
$ws.Cells.Item(9, 2).Value = 400
$ws.Cells.Item(9, 3).Value = 800
$ws.Cells.Item(9, 4).Value = 310
$ws.Cells.Item(9, 5).Value = -400
$ws.Cells.Item(9, 6).Value = -0.5
